{"js": "// ---------------------------------------------------------------------------\n// What the diff actually shows\n// ---------------------------------------------------------------------------\n// The unified diff for this fixture only ever touches attribute ORDER inside\n// existing tags (e.g. <w:pgSz w:w=\"..\" w:h=\"..\"/> -> <w:pgSz w:h=\"..\" w:w=\"..\"/>,\n// the xmlns:* declarations on <w:document>, the <w:rFonts>/<w:lang> attributes\n// in docDefaults, every <w:lsdException>/<w:style> attribute list, etc.) and\n// drops the ephemeral w:rsid* \"who/when touched this run\" bookkeeping\n// attributes that Word stamps on <w:p>/<w:r>/<w:sectPr> while editing. Every\n// \"-\"/\"+\" pair names the same element with the same attribute VALUES; only\n// the serialized order (and the rsid noise) changed. No text, run, paragraph,\n// table, style definition, or numeric value was added, removed, or edited.\n//\n// That lines up exactly with the commit message: \"Fixed POI packaging and\n// upgraded to POI 3.15\" - i.e. the fixture .docx was simply re-saved by a\n// newer revision of the authoring library, which happens to sort XML\n// attributes and drop rsids when it writes the package back out. That is a\n// side effect of the *serializer*, not an editorial change to the document.\n//\n// The Word JavaScript API (like the Word COM object model) edits the\n// document's content/object model - paragraphs, runs, styles, properties -\n// and has no hook for dictating the byte-level attribute order the host\n// chooses when it serializes the underlying OOXML package back to .docx, nor\n// for micromanaging internal rsid bookkeeping. There is therefore no\n// content-level change for this script to make; we simply touch the body\n// read-only to confirm the object model loads/syncs correctly, and leave the\n// document exactly as authored.\n// ---------------------------------------------------------------------------\n\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# ---------------------------------------------------------------------------\n# What the diff actually shows\n# ---------------------------------------------------------------------------\n# The unified diff for this fixture only ever touches attribute ORDER inside\n# existing tags (e.g. <w:pgSz w:w=\"..\" w:h=\"..\"/> -> <w:pgSz w:h=\"..\" w:w=\"..\"/>,\n# the xmlns:* declarations on <w:document>, the <w:rFonts>/<w:lang> attributes\n# in docDefaults, every <w:lsdException>/<w:style> attribute list, etc.) and\n# drops the ephemeral w:rsid* \"who/when touched this run\" bookkeeping\n# attributes that Word stamps on <w:p>/<w:r>/<w:sectPr> while editing. Every\n# \"-\"/\"+\" pair names the same element with the same attribute VALUES; only the\n# serialized order (and the rsid noise) changed. No text, run, paragraph,\n# table, style definition, or numeric value was added, removed, or edited.\n#\n# That lines up exactly with the commit message: \"Fixed POI packaging and\n# upgraded to POI 3.15\" - i.e. the fixture .docx was simply re-saved by a\n# newer revision of the authoring library, which happens to sort XML\n# attributes and drop rsids when it writes the package back out. That is a\n# side effect of the *serializer*, not an editorial change to the document.\n#\n# The Word COM object model (like the Word JavaScript API) edits the\n# document's content/object model - paragraphs, runs, styles, properties -\n# and has no hook for dictating the byte-level attribute order the host\n# chooses when it serializes the underlying OOXML package back to .docx, nor\n# for micromanaging internal rsid bookkeeping. There is therefore no\n# content-level change for this script to make; we simply touch the range\n# read-only to confirm the object model resolves correctly, and leave the\n# document exactly as authored.\n# ---------------------------------------------------------------------------\n\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
